# Update NATMI ligand-receptor pair stats with new TPM-derived values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 17.87134533333333
$ws.Range("H2").Value = 53.614036
$ws.Range("I2").Value = 0.1500697615111392
$ws.Range("J2").Value = 0.1500697615111392
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 2603.868323261825
$ws.Range("R2").Value = 23434.81490935643
$ws.Range("S2").Value = 0.04300905946641546
$ws.Range("T2").Value = 0.04300905946641546
$ws.Range("G3").Value = 17.87134533333333
$ws.Range("H3").Value = 53.614036
$ws.Range("I3").Value = 0.1500697615111392
$ws.Range("J3").Value = 0.1500697615111392
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 3016.677778519988
$ws.Range("R3").Value = 27150.10000667989
$ws.Range("S3").Value = 0.04982758644448325
$ws.Range("T3").Value = 0.04982758644448325
$ws.Range("G4").Value = 17.87134533333333
$ws.Range("H4").Value = 53.614036
$ws.Range("I4").Value = 0.1500697615111392
$ws.Range("J4").Value = 0.1500697615111392
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 2289.786756280078
$ws.Range("R4").Value = 20608.08080652071
$ws.Range("S4").Value = 0.03782125765979369
$ws.Range("T4").Value = 0.03782125765979369
$ws.Range("G5").Value = 17.87134533333333
$ws.Range("H5").Value = 53.614036
$ws.Range("I5").Value = 0.1500697615111392
$ws.Range("J5").Value = 0.1500697615111392
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 1175.238952301615
$ws.Range("R5").Value = 10577.15057071453
$ws.Range("S5").Value = 0.01941185794044679
$ws.Range("T5").Value = 0.01941185794044679
$ws.Range("I6").Value = 0.2793179663930228
$ws.Range("J6").Value = 0.2793179663930228
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 4846.4607225668
$ws.Range("R6").Value = 43618.1465031012
$ws.Range("S6").Value = 0.08005079041685589
$ws.Range("T6").Value = 0.08005079041685589
$ws.Range("I7").Value = 0.2793179663930228
$ws.Range("J7").Value = 0.2793179663930228
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("S7").Value = 0.09274180205125827
$ws.Range("T7").Value = 0.09274180205125827
$ws.Range("I8").Value = 0.2793179663930228
$ws.Range("J8").Value = 0.2793179663930228
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 4261.875102602557
$ws.Range("R8").Value = 38356.87592342301
$ws.Range("S8").Value = 0.07039497277521806
$ws.Range("T8").Value = 0.07039497277521806
$ws.Range("I9").Value = 0.2793179663930228
$ws.Range("J9").Value = 0.2793179663930228
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 2187.418377141804
$ws.Range("R9").Value = 19686.76539427624
$ws.Range("S9").Value = 0.03613040114969055
$ws.Range("T9").Value = 0.03613040114969055
$ws.Range("G10").Value = 12.60542466666667
$ws.Range("H10").Value = 37.816274
$ws.Range("I10").Value = 0.1058506250195358
$ws.Range("J10").Value = 0.1058506250195358
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 1836.619760772902
$ws.Range("R10").Value = 16529.57784695612
$ws.Range("S10").Value = 0.03033613021157856
$ws.Range("T10").Value = 0.03033613021157857
$ws.Range("G11").Value = 12.60542466666667
$ws.Range("H11").Value = 37.816274
$ws.Range("I11").Value = 0.1058506250195358
$ws.Range("J11").Value = 0.1058506250195358
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 2127.791935720399
$ws.Range("R11").Value = 19150.12742148359
$ws.Range("S11").Value = 0.03514552908763041
$ws.Range("T11").Value = 0.03514552908763042
$ws.Range("G12").Value = 12.60542466666667
$ws.Range("H12").Value = 37.816274
$ws.Range("I12").Value = 0.1058506250195358
$ws.Range("J12").Value = 0.1058506250195358
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 1615.084590480348
$ws.Range("R12").Value = 14535.76131432314
$ws.Range("S12").Value = 0.02667695158572574
$ws.Range("T12").Value = 0.02667695158572574
$ws.Range("G13").Value = 12.60542466666667
$ws.Range("H13").Value = 37.816274
$ws.Range("I13").Value = 0.1058506250195358
$ws.Range("J13").Value = 0.1058506250195358
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 828.9463273332152
$ws.Range("R13").Value = 7460.516945998937
$ws.Range("S13").Value = 0.01369201413460109
$ws.Range("T13").Value = 0.01369201413460109
$ws.Range("G14").Value = 55.34703199999999
$ws.Range("H14").Value = 166.041096
$ws.Range("I14").Value = 0.4647616470763022
$ws.Range("J14").Value = 0.4647616470763023
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 8064.103777489829
$ws.Range("R14").Value = 72576.93399740847
$ws.Range("S14").Value = 0.1331977949157343
$ws.Range("T14").Value = 0.1331977949157343
$ws.Range("G15").Value = 55.34703199999999
$ws.Range("H15").Value = 166.041096
$ws.Range("I15").Value = 0.4647616470763022
$ws.Range("J15").Value = 0.4647616470763023
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 9342.562545082485
$ws.Range("R15").Value = 84083.06290574236
$ws.Range("S15").Value = 0.1543145728532122
$ws.Range("T15").Value = 0.1543145728532122
$ws.Range("G16").Value = 55.34703199999999
$ws.Range("H16").Value = 166.041096
$ws.Range("I16").Value = 0.4647616470763022
$ws.Range("J16").Value = 0.4647616470763023
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 7091.402382372949
$ws.Range("R16").Value = 63822.62144135654
$ws.Range("S16").Value = 0.117131324974873
$ws.Range("T16").Value = 0.117131324974873
$ws.Range("G17").Value = 55.34703199999999
$ws.Range("H17").Value = 166.041096
$ws.Range("I17").Value = 0.4647616470763022
$ws.Range("J17").Value = 0.4647616470763023
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 3639.680543767527
$ws.Range("R17").Value = 32757.12489390774
$ws.Range("S17").Value = 0.06011795433248279
$ws.Range("T17").Value = 0.06011795433248279
